$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, new value
$updates = @(
    @(2,  "D", "309.50"),
    @(2,  "E", "-2.74%"),
    @(3,  "D", "37.77"),
    @(3,  "E", "-4.60%"),
    @(4,  "D", "5.091"),
    @(4,  "E", "-1.13%"),
    @(5,  "D", "0.07854"),
    @(5,  "E", "-4.32%"),
    @(6,  "D", "1.958"),
    @(6,  "E", "-6.48%"),
    @(7,  "E", "1.83%"),
    @(8,  "D", "8.300"),
    @(8,  "E", "-0.27%"),
    @(9,  "E", "-6.48%"),
    @(10, "D", "0.9288"),
    @(10, "E", "-0.32%"),
    @(11, "D", "0.1347"),
    @(11, "E", "-3.43%"),
    @(12, "E", "-1.12%"),
    @(13, "D", "0.08983"),
    @(13, "E", "-1.21%"),
    @(14, "D", "0.03472"),
    @(14, "E", "-0.11%"),
    @(15, "E", "-0.95%"),
    @(16, "D", "0.001391"),
    @(16, "E", "-0.47%"),
    @(17, "D", "0.005943"),
    @(17, "E", "-5.85%"),
    @(18, "E", "1,777.07%"),
    @(21, "E", "0.23%"),
    @(22, "D", "5.010"),
    @(22, "E", "3.01%"),
    @(23, "D", "0.2515"),
    @(24, "D", "0.04346"),
    @(24, "E", "0.49%"),
    @(25, "D", "0.001222"),
    @(25, "E", "-0.34%"),
    @(26, "D", "0.004540"),
    @(26, "E", "-4.68%"),
    @(27, "D", "0.0001352"),
    @(27, "E", "3.99%"),
    @(39, "D", "0.02293"),
    @(39, "E", "2.91%"),
    @(40, "D", "0.05049"),
    @(40, "E", "-3.41%"),
    @(41, "D", "0.007611"),
    @(41, "E", "1.36%"),
    @(42, "D", "0.009861"),
    @(42, "E", "2.47%"),
    @(43, "D", "0.1355"),
    @(43, "E", "-1.86%"),
    @(44, "D", "0.002043"),
    @(44, "E", "-4.98%"),
    @(45, "D", "0.008795"),
    @(45, "E", "-10.66%"),
    @(46, "D", "0.00006834"),
    @(46, "E", "3.86%"),
    @(47, "E", "0.17%"),
    @(48, "D", "0.003005"),
    @(48, "E", "8.64%"),
    @(49, "D", "0.001302"),
    @(49, "E", "8.50%"),
    @(50, "D", "0.00002103"),
    @(50, "E", "0.17%"),
    @(51, "D", "0.0002003"),
    @(51, "E", "0.17%")
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $val = $u[2]
    $cell = $ws.Range("$col$row")
    # Force text interpretation so values like "309.50" or "-2.74%" are
    # stored as literal strings (matching the original inlineStr cells)
    # instead of being parsed into numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $val
}
